# 15 January 2024 - Added random data points to the two new variables.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Age (column D) and Political Party (column E) values for rows 2-15
$ages = @(17, 25, 89, 65, 49, 33, 23, 9, 12, 99, 26, 56, 77, 45)
$parties = @("Republican", "Democrat", "Democrat", "Democrat", "Republican", "Libertarian", "Green", "Libertarian", "Republican", "Democrat", "Libertarian", "Democrat", "Democrat", "Democrat")

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $ages[$i]
    $ws.Cells.Item($row, 5).Value = $parties[$i]
}

# Update the active selection to match the final edited cell (E16)
$ws.Range("E16").Select()
